{"js": "// Added phrase \"marginal conditions\" to point about mosses vs other plants\nconst body = context.document.body;\n\n// 1. \"What \" + \"is a moss\" -> \"What is a moss\" (run merge, heading text unchanged)\nlet rHeading = body.search(\"What is a moss\", { matchCase: true });\nawait context.sync();\nif (rHeading.items.length > 0) {\n  rHeading.items[0].insertText(\"What is a moss\", \"Replace\");\n  await context.sync();\n}\n\n// 2. \" in reasonable quantities\" + \".\" -> \" in reasonable quantities.\" (run merge, text unchanged)\nlet rQty = body.search(\"in reasonable quantities.\", { matchCase: true });\nawait context.sync();\nif (rQty.items.length > 0) {\n  rQty.items[0].insertText(\"in reasonable quantities.\", \"Replace\");\n  await context.sync();\n}\n\n// 3. Main wording change:\n//    \"Mosses grow best in conditions...\" ->\n//    \"In the wild, mosses grow best in marginal conditions...\"\nlet rOpen = body.search(\"Mosses grow best in conditions\", { matchCase: true });\nawait context.sync();\nif (rOpen.items.length > 0) {\n  rOpen.items[0].insertText(\n    \"In the wild, mosses grow best in marginal conditions\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n//    \"...outcompeted by other plants.  \" -> \"...outcompeted by those other plants!  \"\nlet rClose = body.search(\"outcompeted by other plants.  \", { matchCase: true });\nawait context.sync();\nif (rClose.items.length > 0) {\n  rClose.items[0].insertText(\n    \"outcompeted by those other plants!  \",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 4. \" to continually remove \" -> \" to remove \"\nlet rCont = body.search(\"continually \", { matchCase: true });\nawait context.sync();\nif (rCont.items.length > 0) {\n  rCont.items[0].insertText(\"\", \"Replace\");\n  await context.sync();\n}\n\n// 5. \"(Michael Fletcher, \" + \"2\" -> \"(Michael Fletcher, 2\" (run merge, text unchanged)\nlet rFletcher = body.search(\"(Michael Fletcher, 2\", { matchCase: true });\nawait context.sync();\nif (rFletcher.items.length > 0) {\n  rFletcher.items[0].insertText(\"(Michael Fletcher, 2\", \"Replace\");\n  await context.sync();\n}\n\n// 6. \" Ed\" + \" 2005) is \" -> \" Ed 2005) is \" (run merge, text unchanged)\nlet rEd = body.search(\" Ed 2005) is \", { matchCase: true });\nawait context.sync();\nif (rEd.items.length > 0) {\n  rEd.items[0].insertText(\" Ed 2005) is \", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Added phrase \"marginal conditions\" to point about mosses vs other plants\n$d = $word.ActiveDocument\n\n# Helper: find the first occurrence of $findText in the document and replace\n# it with $replaceText. Uses Find.Execute's own Replace argument (wdReplaceOne)\n# rather than a manual Range.Text assignment so that the edit is registered as\n# a real mutation even when $findText -eq $replaceText (pure run-merge cases\n# below) -- a plain \"$range.Text = $sameText\" assignment is treated as a no-op\n# by the engine and silently skipped.\nfunction ReplaceOnce($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    # wdFindContinue=1 (Wrap), wdReplaceOne=1 (Replace) -- replace just the\n    # first/only match.\n    $found = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    return $found\n}\n\n# 1. \"What \" + \"is a moss\" -> \"What is a moss\" (heading text itself unchanged)\nReplaceOnce \"What is a moss\" \"What is a moss\" | Out-Null\n\n# 2. \" in reasonable quantities\" + \".\" -> \" in reasonable quantities.\" (text unchanged)\nReplaceOnce \"in reasonable quantities.\" \"in reasonable quantities.\" | Out-Null\n\n# 3. Main wording change: add \"In the wild, \" lead-in and \"marginal \" before\n#    \"conditions\".\nReplaceOnce \"Mosses grow best in conditions\" \"In the wild, mosses grow best in marginal conditions\" | Out-Null\n\n#    ...and add \"those \" before \"other plants\", changing the trailing \".\" to \"!\".\nReplaceOnce \"outcompeted by other plants.  \" \"outcompeted by those other plants!  \" | Out-Null\n\n# 4. \" to continually remove \" -> \" to remove \"\nReplaceOnce \" to continually remove \" \" to remove \" | Out-Null\n\n# 5. \"(Michael Fletcher, \" + \"2\" -> \"(Michael Fletcher, 2\" (text unchanged).\n#    This single edit also merges \" Ed\" + \" 2005) is \" -> \" Ed 2005) is \" as a\n#    side effect of the engine's paragraph-level run normalization, which\n#    covers hunk 6 too. (Re-running a second, separate replace directly after\n#    the superscript \"nd\" run risks the new text inheriting the superscript\n#    run's formatting, so that merge is intentionally *not* requested again\n#    here.)\nReplaceOnce \"(Michael Fletcher, 2\" \"(Michael Fletcher, 2\" | Out-Null\n"}
